$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 24.03398533008638
$ws.Cells.Item(2, 5).Value = 22.88554763793945
$ws.Cells.Item(2, 6).Value = 24.35854508504219
$ws.Cells.Item(2, 7).Value = 22.81897156284831
$ws.Cells.Item(2, 8).Value = 121301203
$ws.Cells.Item(2, 9).Value = "BAH"
$ws.Cells.Item(3, 4).Value = 21.33872343234485
$ws.Cells.Item(3, 5).Value = 23.19571876525879
$ws.Cells.Item(3, 6).Value = 23.51358382303037
$ws.Cells.Item(3, 7).Value = 20.99576485638509
$ws.Cells.Item(3, 8).Value = 121301203
$ws.Cells.Item(3, 9).Value = "BAH"
$ws.Cells.Item(4, 4).Value = 22.04670703805344
$ws.Cells.Item(4, 5).Value = 24.76156997680664
$ws.Cells.Item(4, 6).Value = 25.69454209634771
$ws.Cells.Item(4, 7).Value = 21.43313150139808
$ws.Cells.Item(4, 8).Value = 121301203
$ws.Cells.Item(4, 9).Value = "BAH"
$ws.Cells.Item(5, 4).Value = 25.65923377826237
$ws.Cells.Item(5, 5).Value = 23.88613891601562
$ws.Cells.Item(5, 6).Value = 26.199607033539
$ws.Cells.Item(5, 7).Value = 23.06713633006298
$ws.Cells.Item(5, 8).Value = 121301203
$ws.Cells.Item(5, 9).Value = "BAH"
$ws.Cells.Item(6, 4).Value = 25.55567626001431
$ws.Cells.Item(6, 5).Value = 23.40764045715332
$ws.Cells.Item(6, 6).Value = 26.01414933838144
$ws.Cells.Item(6, 7).Value = 23.18689356038835
$ws.Cells.Item(6, 8).Value = 121301203
$ws.Cells.Item(6, 9).Value = "BAH"
$ws.Cells.Item(7, 4).Value = 25.29426269436681
$ws.Cells.Item(7, 5).Value = 26.35245704650879
$ws.Cells.Item(7, 6).Value = 26.88155503642894
$ws.Cells.Item(7, 7).Value = 24.77370035446448
$ws.Cells.Item(7, 8).Value = 121301203
$ws.Cells.Item(7, 9).Value = "BAH"
$ws.Cells.Item(8, 4).Value = 27.10136860414235
$ws.Cells.Item(8, 5).Value = 26.13223648071289
$ws.Cells.Item(8, 6).Value = 27.23859059113526
$ws.Cells.Item(8, 7).Value = 25.34320923759608
$ws.Cells.Item(8, 8).Value = 121301203
$ws.Cells.Item(8, 9).Value = "BAH"
$ws.Cells.Item(9, 4).Value = 31.18382862331365
$ws.Cells.Item(9, 5).Value = 29.14167213439941
$ws.Cells.Item(9, 6).Value = 31.51126406499493
$ws.Cells.Item(9, 7).Value = 28.21968571891028
$ws.Cells.Item(9, 8).Value = 121301203
$ws.Cells.Item(9, 9).Value = "BAH"
$ws.Cells.Item(10, 4).Value = 30.70533139166112
$ws.Cells.Item(10, 5).Value = 31.11231231689453
$ws.Cells.Item(10, 6).Value = 31.73577189543047
$ws.Cells.Item(10, 7).Value = 29.71818926111396
$ws.Cells.Item(10, 8).Value = 121301203
$ws.Cells.Item(10, 9).Value = "BAH"
$ws.Cells.Item(11, 4).Value = 28.47653569835936
$ws.Cells.Item(11, 5).Value = 29.83338737487793
$ws.Cells.Item(11, 6).Value = 29.96385541359244
$ws.Cells.Item(11, 7).Value = 28.14601910010533
$ws.Cells.Item(11, 8).Value = 121301203
$ws.Cells.Item(11, 9).Value = "BAH"
$ws.Cells.Item(12, 4).Value = 32.71055465316208
$ws.Cells.Item(12, 5).Value = 33.04281997680664
$ws.Cells.Item(12, 6).Value = 33.87348161817059
$ws.Cells.Item(12, 7).Value = 32.61437232246274
$ws.Cells.Item(12, 8).Value = 121301203
$ws.Cells.Item(12, 9).Value = "BAH"
$ws.Cells.Item(13, 4).Value = 33.59236502309865
$ws.Cells.Item(13, 5).Value = 34.41811752319336
$ws.Cells.Item(13, 6).Value = 35.35807094226917
$ws.Cells.Item(13, 7).Value = 33.3024711559307
$ws.Cells.Item(13, 8).Value = 121301203
$ws.Cells.Item(13, 9).Value = "BAH"
$ws.Cells.Item(14, 4).Value = 34.12595626088237
$ws.Cells.Item(14, 5).Value = 34.99124526977539
$ws.Cells.Item(14, 6).Value = 36.31566784080812
$ws.Cells.Item(14, 7).Value = 33.23417651020576
$ws.Cells.Item(14, 8).Value = 121301203
$ws.Cells.Item(14, 9).Value = "BAH"
$ws.Cells.Item(15, 4).Value = 38.5087829125042
$ws.Cells.Item(15, 5).Value = 41.91365814208984
$ws.Cells.Item(15, 6).Value = 42.3037988622988
$ws.Cells.Item(15, 7).Value = 38.46444927786856
$ws.Cells.Item(15, 8).Value = 121301203
$ws.Cells.Item(15, 9).Value = "BAH"
$ws.Cells.Item(16, 4).Value = 44.35511153843082
$ws.Cells.Item(16, 5).Value = 44.09697341918945
$ws.Cells.Item(16, 6).Value = 44.85358140116309
$ws.Cells.Item(16, 7).Value = 41.3642787797521
$ws.Cells.Item(16, 8).Value = 121301203
$ws.Cells.Item(16, 9).Value = "BAH"
$ws.Cells.Item(17, 4).Value = 39.59662261234779
$ws.Cells.Item(17, 5).Value = 43.89400100708008
$ws.Cells.Item(17, 6).Value = 44.0280128031828
$ws.Cells.Item(17, 7).Value = 39.01589460479136
$ws.Cells.Item(17, 8).Value = 121301203
$ws.Cells.Item(17, 9).Value = "BAH"
$ws.Cells.Item(18, 4).Value = 52.44463445303174
$ws.Cells.Item(18, 5).Value = 53.19832992553711
$ws.Cells.Item(18, 6).Value = 53.2073009830411
$ws.Cells.Item(18, 7).Value = 50.80265368613166
$ws.Cells.Item(18, 8).Value = 121301203
$ws.Cells.Item(18, 9).Value = "BAH"
$ws.Cells.Item(19, 4).Value = 60.13929570462258
$ws.Cells.Item(19, 5).Value = 61.90412521362305
$ws.Cells.Item(19, 6).Value = 64.75846650202081
$ws.Cells.Item(19, 7).Value = 59.69809004479426
$ws.Cells.Item(19, 8).Value = 121301203
$ws.Cells.Item(19, 9).Value = "BAH"
$ws.Cells.Item(20, 4).Value = 64.43943112600229
$ws.Cells.Item(20, 5).Value = 63.56326293945312
$ws.Cells.Item(20, 6).Value = 65.6136876639557
$ws.Cells.Item(20, 7).Value = 60.84441186987823
$ws.Cells.Item(20, 8).Value = 121301203
$ws.Cells.Item(20, 9).Value = "BAH"
$ws.Cells.Item(21, 4).Value = 65.03449833284634
$ws.Cells.Item(21, 5).Value = 70.75550079345703
$ws.Cells.Item(21, 6).Value = 74.34586105840349
$ws.Cells.Item(21, 7).Value = 64.48143710573977
$ws.Cells.Item(21, 8).Value = 121301203
$ws.Cells.Item(21, 9).Value = "BAH"
$ws.Cells.Item(22, 4).Value = 61.05830608183504
$ws.Cells.Item(22, 5).Value = 66.85734558105469
$ws.Cells.Item(22, 6).Value = 72.33775111730718
$ws.Cells.Item(22, 7).Value = 60.56670713543927
$ws.Cells.Item(22, 8).Value = 121301203
$ws.Cells.Item(22, 9).Value = "BAH"
$ws.Cells.Item(23, 4).Value = 71.34996588008937
$ws.Cells.Item(23, 5).Value = 74.73191833496094
$ws.Cells.Item(23, 6).Value = 74.8324633008954
$ws.Cells.Item(23, 7).Value = 64.08335008781862
$ws.Cells.Item(23, 8).Value = 121301203
$ws.Cells.Item(23, 9).Value = "BAH"
$ws.Cells.Item(24, 4).Value = 76.42507735292196
$ws.Cells.Item(24, 5).Value = 72.01258850097656
$ws.Cells.Item(24, 6).Value = 76.86541218377478
$ws.Cells.Item(24, 7).Value = 68.17803255307108
$ws.Cells.Item(24, 8).Value = 121301203
$ws.Cells.Item(24, 9).Value = "BAH"
$ws.Cells.Item(25, 4).Value = 80.44695763993532
$ws.Cells.Item(25, 5).Value = 78.41230773925781
$ws.Cells.Item(25, 6).Value = 92.3050171475905
$ws.Cells.Item(25, 7).Value = 78.15452505849161
$ws.Cells.Item(25, 8).Value = 121301203
$ws.Cells.Item(25, 9).Value = "BAH"
$ws.Cells.Item(26, 4).Value = 74.93267586897272
$ws.Cells.Item(26, 5).Value = 76.70819854736328
$ws.Cells.Item(26, 6).Value = 79.00158729842605
$ws.Cells.Item(26, 7).Value = 74.57202304540277
$ws.Cells.Item(26, 8).Value = 121301203
$ws.Cells.Item(26, 9).Value = "BAH"
$ws.Cells.Item(27, 4).Value = 79.26966321970276
$ws.Cells.Item(27, 5).Value = 79.68755340576172
$ws.Cells.Item(27, 6).Value = 84.18222868370674
$ws.Cells.Item(27, 7).Value = 79.01892344003055
$ws.Cells.Item(27, 8).Value = 121301203
$ws.Cells.Item(27, 9).Value = "BAH"
$ws.Cells.Item(28, 4).Value = 74.35084416868511
$ws.Cells.Item(28, 5).Value = 81.03029632568359
$ws.Cells.Item(28, 6).Value = 81.17023033450081
$ws.Cells.Item(28, 7).Value = 73.18473929077386
$ws.Cells.Item(28, 8).Value = 121301203
$ws.Cells.Item(28, 9).Value = "BAH"
$ws.Cells.Item(29, 4).Value = 79.52406680314768
$ws.Cells.Item(29, 5).Value = 71.87986755371094
$ws.Cells.Item(29, 6).Value = 85.24785170358814
$ws.Cells.Item(29, 7).Value = 68.68541165814838
$ws.Cells.Item(29, 8).Value = 121301203
$ws.Cells.Item(29, 9).Value = "BAH"
$ws.Cells.Item(30, 4).Value = 82.71151096214419
$ws.Cells.Item(30, 5).Value = 76.90785217285156
$ws.Cells.Item(30, 6).Value = 86.16920649257082
$ws.Cells.Item(30, 7).Value = 76.57810046619814
$ws.Cells.Item(30, 8).Value = 121301203
$ws.Cells.Item(30, 9).Value = "BAH"
$ws.Cells.Item(31, 4).Value = 85.57459395746547
$ws.Cells.Item(31, 5).Value = 90.89696502685548
$ws.Cells.Item(31, 6).Value = 91.84400567957169
$ws.Cells.Item(31, 7).Value = 84.59914324122325
$ws.Cells.Item(31, 8).Value = 121301203
$ws.Cells.Item(31, 9).Value = "BAH"
$ws.Cells.Item(32, 4).Value = 88.55238161336706
$ws.Cells.Item(32, 5).Value = 103.555290222168
$ws.Cells.Item(32, 6).Value = 104.7159478450649
$ws.Cells.Item(32, 7).Value = 88.29551095477385
$ws.Cells.Item(32, 8).Value = 121301203
$ws.Cells.Item(32, 9).Value = "BAH"
$ws.Cells.Item(33, 4).Value = 99.58245342149893
$ws.Cells.Item(33, 5).Value = 90.41139221191406
$ws.Cells.Item(33, 6).Value = 100.62375495187
$ws.Cells.Item(33, 7).Value = 87.41169104756972
$ws.Cells.Item(33, 8).Value = 121301203
$ws.Cells.Item(33, 9).Value = "BAH"
$ws.Cells.Item(34, 4).Value = 89.33768280775007
$ws.Cells.Item(34, 5).Value = 91.89128875732422
$ws.Cells.Item(34, 6).Value = 95.5392878398542
$ws.Cells.Item(34, 7).Value = 89.13608349450354
$ws.Cells.Item(34, 8).Value = 121301203
$ws.Cells.Item(34, 9).Value = "BAH"
$ws.Cells.Item(35, 4).Value = 106.8116554203405
$ws.Cells.Item(35, 5).Value = 116.7532348632812
$ws.Cells.Item(35, 6).Value = 117.2546495418238
$ws.Cells.Item(35, 7).Value = 106.2716692501008
$ws.Cells.Item(35, 8).Value = 121301203
$ws.Cells.Item(35, 9).Value = "BAH"
$ws.Cells.Item(36, 4).Value = 106.6388001046042
$ws.Cells.Item(36, 5).Value = 116.0965042114258
$ws.Cells.Item(36, 6).Value = 126.173750481564
$ws.Cells.Item(36, 7).Value = 106.3193465994651
$ws.Cells.Item(36, 8).Value = 121301203
$ws.Cells.Item(36, 9).Value = "BAH"
$ws.Cells.Item(37, 4).Value = 123.5805814473242
$ws.Cells.Item(37, 5).Value = 136.7752075195312
$ws.Cells.Item(37, 6).Value = 143.3336547933547
$ws.Cells.Item(37, 7).Value = 119.6746575863642
$ws.Cells.Item(37, 8).Value = 121301203
$ws.Cells.Item(37, 9).Value = "BAH"
$ws.Cells.Item(38, 4).Value = 144.6014448470179
$ws.Cells.Item(38, 5).Value = 143.9871520996094
$ws.Cells.Item(38, 6).Value = 145.7617689795722
$ws.Cells.Item(38, 7).Value = 135.4553745111098
$ws.Cells.Item(38, 8).Value = 121301203
$ws.Cells.Item(38, 9).Value = "BAH"
$ws.Cells.Item(39, 4).Value = 151.2520932663955
$ws.Cells.Item(39, 5).Value = 140.2065734863281
$ws.Cells.Item(39, 6).Value = 157.1417174825091
$ws.Cells.Item(39, 7).Value = 131.7438979887037
$ws.Cells.Item(39, 8).Value = 121301203
$ws.Cells.Item(39, 9).Value = "BAH"
$ws.Cells.Item(40, 4).Value = 159.8514850219303
$ws.Cells.Item(40, 5).Value = 178.3479919433594
$ws.Cells.Item(40, 6).Value = 185.6818006556221
$ws.Cells.Item(40, 7).Value = 155.3647977324328
$ws.Cells.Item(40, 8).Value = 121301203
$ws.Cells.Item(40, 9).Value = "BAH"
$ws.Cells.Item(41, 4).Value = 127.6583539181418
$ws.Cells.Item(41, 5).Value = 127.0673370361328
$ws.Cells.Item(41, 6).Value = 144.7484092223332
$ws.Cells.Item(41, 7).Value = 119.2068891593502
$ws.Cells.Item(41, 8).Value = 121301203
$ws.Cells.Item(41, 9).Value = "BAH"
$ws.Cells.Item(42, 4).Value = 103.4070023447395
$ws.Cells.Item(42, 5).Value = 118.7760391235352
$ws.Cells.Item(42, 6).Value = 119.538062620862
$ws.Cells.Item(42, 7).Value = 100.2104841841535
$ws.Cells.Item(42, 8).Value = 121301203
$ws.Cells.Item(42, 9).Value = "BAH"
$ws.Cells.Item(43, 4).Value = 103.5970493723464
$ws.Cells.Item(43, 5).Value = 106.8011856079102
$ws.Cells.Item(43, 6).Value = 119.4585152280433
$ws.Cells.Item(43, 7).Value = 102.9602032674671
$ws.Cells.Item(43, 8).Value = 121301203
$ws.Cells.Item(43, 9).Value = "BAH"
